$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column L ("correct_ans") used abbreviated codes (r/y/b); spell them out ---
$ansMap = @{ "r" = "right"; "y" = "left"; "b" = "center" }

$usedRange  = $ws.UsedRange
$lastRow    = $usedRange.Rows.Count
$lastDataRow = $lastRow  # header is row 1, data starts row 2

$colL = 12  # column L = correct_ans

for ($r = 2; $r -le $lastDataRow; $r++) {
    $cell = $ws.Cells.Item($r, $colL)
    $old = $cell.Value2
    if ($null -ne $old -and $ansMap.ContainsKey($old)) {
        $cell.Value2 = $ansMap[$old]
    }
}

# --- 2) Stimulus set renamed: "face//face_NN.jpg" -> "book//book_NN.jpg" ---
# These filenames only ever show up in columns A-D (promptFile/correctFile/dist_01File/dist_02File).
$firstFileCol = 1  # column A
$lastFileCol  = 4  # column D

for ($r = 2; $r -le $lastDataRow; $r++) {
    for ($c = $firstFileCol; $c -le $lastFileCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val -and $val -like "*face//face_*") {
            $cell.Value2 = $val.Replace("face//face_", "book//book_")
        }
    }
}
